$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 478
$ws.Range("E4").Value = 3
$ws.Range("I4").Value = 0.01255230125523013
$ws.Range("J4").Value = 0.3
$ws.Range("K4").Value = 0.02409638554216867
$ws.Range("L4").Value = 0.006276150627615063
$ws.Range("M4").Value = 0.375
$ws.Range("N4").Value = 0.01234567901234568
